$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 76
$ws.Range("H76").Value = 100006530
$ws.Range("I76").Value = 137505970
$ws.Range("J76").Value = 8000
$ws.Range("K76").Value = 137505970
$ws.Range("L76").Value = 8000
$ws.Range("M76").Value = -137505655
$ws.Range("N76").Value = -8630

# Row 79
$ws.Range("H79").Value = 100006530
$ws.Range("I79").Value = 137505970
$ws.Range("J79").Value = 8000
$ws.Range("K79").Value = 137505970
$ws.Range("L79").Value = 8000
$ws.Range("M79").Value = -137504878
$ws.Range("N79").Value = -10184

# Row 113
$ws.Range("H113").Value = 2393.5417
$ws.Range("I113").Value = 2617
$ws.Range("J113").Value = 2021.1111
$ws.Range("K113").Value = 2617
$ws.Range("L113").Value = 2021.1111
$ws.Range("M113").Value = 637
$ws.Range("N113").Value = -8529.1111

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6322.8203
$ws.Range("I32").Value = 5268.9697
$ws.Range("K32").Value = 5268.9697
$ws.Range("M32").Value = -4981.9697

# Row 63
$ws.Range("H63").Value = 1754
$ws.Range("I63").Value = 1742.5
$ws.Range("K63").Value = 1742.5
$ws.Range("M63").Value = -1056.5

# Row 66
$ws.Range("H66").Value = 1754
$ws.Range("I66").Value = 1742.5
$ws.Range("K66").Value = 8712.5
$ws.Range("M66").Value = -5280.5

# Row 74
$ws.Range("H74").Value = 1066.0476
$ws.Range("I74").Value = 1007.25
$ws.Range("J74").Value = 1144.4445
$ws.Range("K74").Value = 1007.25
$ws.Range("L74").Value = 1144.4445
$ws.Range("M74").Value = -133.25
$ws.Range("N74").Value = -2892.4445

# Row 77
$ws.Range("H77").Value = 1066.0476
$ws.Range("I77").Value = 1007.25
$ws.Range("J77").Value = 1144.4445
$ws.Range("K77").Value = 5036.25
$ws.Range("L77").Value = 5722.2225
$ws.Range("M77").Value = -668.25
$ws.Range("N77").Value = -14458.2225

# Row 110
$ws.Range("H110").Value = 964.4
$ws.Range("I110").Value = 964.4
$ws.Range("K110").Value = 964.4
$ws.Range("M110").Value = 1080.6

# Row 132
$ws.Range("H132").Value = 3961.9375
$ws.Range("I132").Value = 3489.4
$ws.Range("J132").Value = 4749.5
$ws.Range("K132").Value = 10468.2
$ws.Range("L132").Value = 14248.5
$ws.Range("M132").Value = -7938.200000000001
$ws.Range("N132").Value = -19308.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3052.7273
$ws.Range("I31").Value = 2968
$ws.Range("K31").Value = 2968
$ws.Range("M31").Value = -2673

# Row 34
$ws.Range("H34").Value = 3052.7273
$ws.Range("I34").Value = 2968
$ws.Range("K34").Value = 2968
$ws.Range("M34").Value = -2766

# Row 82
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

# Row 85
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 63
$ws.Range("H63").Value = 3378.9
$ws.Range("I63").Value = 1757.8
$ws.Range("K63").Value = 5273.4
$ws.Range("M63").Value = -4524.4

# Row 66
$ws.Range("H66").Value = 3378.9
$ws.Range("I66").Value = 1757.8
$ws.Range("K66").Value = 15820.2
$ws.Range("M66").Value = -12076.2

# Row 131
$ws.Range("H131").Value = 13186329
$ws.Range("J131").Value = 2476.6296
$ws.Range("L131").Value = 7429.888800000001
$ws.Range("N131").Value = -17509.8888

$ws = $wb.Worksheets.Item("GSM")
# Row 19
$ws.Range("H19").Value = 225282.44
$ws.Range("I19").Value = 1000150
$ws.Range("J19").Value = 3891.7144
$ws.Range("K19").Value = 1000150
$ws.Range("L19").Value = 3891.7144
$ws.Range("M19").Value = -999862
$ws.Range("N19").Value = -4467.7144

# Row 43
$ws.Range("H43").Value = 8919.833000000001
$ws.Range("J43").Value = 8919.833000000001
$ws.Range("L43").Value = 8919.833000000001
$ws.Range("N43").Value = -9221.833000000001

# Row 46
$ws.Range("H46").Value = 20933.334
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 20933.334
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 20933.334
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -21245.334

# Row 70
$ws.Range("H70").Value = 4469.5
$ws.Range("I70").Value = 4009.8572
$ws.Range("J70").Value = 6400
$ws.Range("K70").Value = 4009.8572
$ws.Range("L70").Value = 6400
$ws.Range("M70").Value = -3739.8572
$ws.Range("N70").Value = -6940

# Row 73
$ws.Range("H73").Value = 4469.5
$ws.Range("I73").Value = 4009.8572
$ws.Range("J73").Value = 6400
$ws.Range("K73").Value = 4009.8572
$ws.Range("L73").Value = 6400
$ws.Range("M73").Value = -3073.8572
$ws.Range("N73").Value = -8272

# Row 102
$ws.Range("H102").Value = 835.75
$ws.Range("I102").Value = 835.75
$ws.Range("K102").Value = 835.75
$ws.Range("M102").Value = 786.25

# Row 107
$ws.Range("H107").Value = 1050940.4
$ws.Range("I107").Value = 2262976
$ws.Range("J107").Value = 509.6
$ws.Range("K107").Value = 2262976
$ws.Range("L107").Value = 509.6
$ws.Range("M107").Value = -2261056
$ws.Range("N107").Value = -4349.6

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1508.2222
$ws.Range("I7").Value = 1506.1875
$ws.Range("J7").Value = 1524.5
$ws.Range("K7").Value = 1506.1875
$ws.Range("L7").Value = 1524.5
$ws.Range("M7").Value = -1394.1875
$ws.Range("N7").Value = -1748.5

# Row 22
$ws.Range("H22").Value = 509.94446
$ws.Range("I22").Value = 525.3333
$ws.Range("J22").Value = 433
$ws.Range("K22").Value = 525.3333
$ws.Range("L22").Value = 433
$ws.Range("M22").Value = -230.3333
$ws.Range("N22").Value = -1023

# Row 27
$ws.Range("H27").Value = 509.94446
$ws.Range("I27").Value = 525.3333
$ws.Range("J27").Value = 433
$ws.Range("K27").Value = 525.3333
$ws.Range("L27").Value = 433
$ws.Range("M27").Value = -418.3333
$ws.Range("N27").Value = -647

# Row 126
$ws.Range("H126").Value = 1508.2222
$ws.Range("I126").Value = 1506.1875
$ws.Range("J126").Value = 1524.5
$ws.Range("K126").Value = 4518.5625
$ws.Range("L126").Value = 4573.5
$ws.Range("M126").Value = -2048.5625
$ws.Range("N126").Value = -9513.5

$ws = $wb.Worksheets.Item("WVR")
# Row 140
$ws.Range("H140").Value = 71280
$ws.Range("J140").Value = 71280
$ws.Range("L140").Value = 71280
